$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '80.804.43'
$ws.Range('E2').Value = '  +1.60%  '
$ws.Range('D3').Value = '3.136.81'
$ws.Range('E3').Value = '  -2.16%  '
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').Value = '''204.94'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.16%  '
$ws.Range('D6').Value = '''622.63'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.27%  '
$ws.Range('D7').Value = '''0.279'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +22.53%  '
$ws.Range('D8').Value = '''1.00'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('D9').Value = '''0.577'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.96%  '
$ws.Range('D10').Value = '3.138.32'
$ws.Range('E10').Value = '  -2.09%  '
$ws.Range('D11').Value = '''0.574'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.67%  '
$ws.Range('D12').Value = '''0.0000249'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +8.87%  '
$ws.Range('E13').Value = '  +0.92%  '
$ws.Range('D14').Value = '''5.26'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.11%  '
$ws.Range('D15').Value = '3.718.11'
$ws.Range('E15').Value = '  -2.02%  '
$ws.Range('D16').Value = '''31.15'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.14%  '
$ws.Range('D17').Value = '80.836.77'
$ws.Range('E17').Value = '  +2.20%  '
$ws.Range('D18').Value = '3.150.25'
$ws.Range('E18').Value = '  -1.57%  '
$ws.Range('D19').Value = '''3.13'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +9.71%  '
$ws.Range('D20').Value = '''13.89'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.17%  '
$ws.Range('D21').Value = '''430.10'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.78%  '
$ws.Range('D22').Value = '''8.93'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.70%  '
$ws.Range('D23').Value = '''5.06'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.47%  '
$ws.Range('D24').Value = '''7.13'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.64%  '
$ws.Range('D25').Value = '''5.15'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +7.87%  '
$ws.Range('D26').Value = '3.304.82'
$ws.Range('E26').Value = '  -1.97%  '
$ws.Range('B27').Value = 'Litecoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D27').Value = '''75.64'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.99%  '
$ws.Range('B28').Value = 'Aptos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D28').Value = '''10.89'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.37%  '
$ws.Range('D29').Value = '''0.999'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.38%  '
$ws.Range('D30').Value = '''0.0000121'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +5.09%  '
$ws.Range('E31').Value = '  +0.31%  '
$ws.Range('D32').Value = '''8.95'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.18%  '
$ws.Range('B33').Value = 'Bittensor'
$ws.Range('C33').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D33').Value = '''554.71'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +6.89%  '
$ws.Range('B34').Value = 'Cronos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D34').Value = '''0.151'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +24.56%  '
$ws.Range('D35').Value = '''1.47'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').Value = '''0.150'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +8.92%  '
$ws.Range('B37').Value = 'PancakeSwap'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D37').Value = '''1.99'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.71%  '
$ws.Range('D38').Value = '''22.60'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.93%  '
$ws.Range('D39').Value = '''1.00'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.03%  '
$ws.Range('D40').Value = '''0.404'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.44%  '
$ws.Range('D41').Value = '''5.89'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +8.10%  '
$ws.Range('D42').Value = '''20.70'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.41%  '
$ws.Range('D43').Value = '''3.01'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +18.48%  '
$ws.Range('B44').Value = 'Monero'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D44').Value = '''160.37'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.40%  '
$ws.Range('B45').Value = 'Stacks'
$ws.Range('C45').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D45').Value = '''1.97'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +10.32%  '
$ws.Range('E46').Value = '  -0.02%  '
$ws.Range('D47').Value = '''185.85'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.76%  '
$ws.Range('D48').Value = '''1.31'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.06%  '
$ws.Range('D49').Value = '''43.79'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.50%  '
$ws.Range('D50').Value = '''0.770'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.18%  '
$ws.Range('B51').Value = 'Filecoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D51').Value = '''4.20'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.71%  '
